# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row number -> new value for column F
$updates = @{
    2  = 154
    3  = 1714
    4  = 790
    7  = 11978
    10 = 477
    13 = 857
    14 = 13473
    15 = 13453
    23 = 294
    24 = 172
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
